$wb = $excel.ActiveWorkbook
$ws9 = $wb.Worksheets.Item("2000 vs test 2000")

# Reorder existing rows 7-10 and append two new rows (11-12).
# Target order (label, error_rate):
#   7  Pool5   0.0095
#   8  Fc6     0.0125
#   9  Fc7     0.0075
#   10 Fc8     0.007
#   11 Prob    0.027
#   12 newFc7  0.007
$labels = @("Pool5", "Fc6", "Fc7", "Fc8", "Prob", "newFc7")
$values = @(0.0095, 0.0125, 0.0075, 0.007, 0.027, 0.007)

for ($i = 0; $i -lt $labels.Length; $i++) {
    $r = 7 + $i
    $ws9.Cells.Item($r, 1).Value = $labels[$i]
    $ws9.Cells.Item($r, 2).Value = $values[$i]
}

# Rows 7-10 reuse the already-present "accuracy" formula cells (edited in
# place, one at a time, so the existing shared-formula group is preserved).
for ($r = 7; $r -le 10; $r++) {
    $ws9.Cells.Item($r, 3).Formula = "=1-B$r"
}

# Rows 11-12 are brand new - fill them as a single range so the two cells
# share one formula group, just like a fill-down of C10 would.
$ws9.Range("C11:C12").Formula = "=1-B11"

# Activate the "2000 vs test 2000" sheet and leave the selection on A12,
# matching the saved view state in the target workbook.
$ws9.Activate()
$ws9.Range("A12").Select() | Out-Null
